# Applies the "Fixed citation graph with full text" edit described by the diff.
# The workbook has a header row (row 1) and 6 data rows (rows 2-7) describing
# citation records pulled from PMC / CrossRef. This script updates the
# Abstract/Authors/"Other found locations"/"Misc. Data" columns (D, E, I, J)
# for rows 2-5, and resets the CrossRef-only rows (6-7) back to "unknown"
# placeholder values across Title/Abstract/Authors/ID/ID Format/Date
# Accepted/Misc. Data (C, D, E, F, G, H, J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Egol et al., Increased Mortality...) ---
$ws.Range("D2").Value = "Supplemental Digital Content is Available in the Text.`n"
$ws.Range("E2").Value = "[Kenneth A.%Egol%NULL%0, Sanjit R.%Konda%NULL%2, Mackenzie L.%Bird%NULL%2, Nicket%Dedhia%NULL%2, Emma K.%Landes%NULL%2, Rachel A.%Ranson%NULL%2, Sara J.%Solasz%NULL%2, Vinay K.%Aggarwal%NULL%2, Joseph A.%Bosco%NULL%2, David L.%Furgiuele%NULL%2, Abhishek%Ganta%NULL%2, Jason%Gould%NULL%2, Thomas R.%Lyon%NULL%2, Toni M.%McLaurin%NULL%2, Nirmal C.%Tejwani%NULL%2, Joseph D.%Zuckerman%NULL%2, Philipp%Leucht%NULL%2]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Journal of Orthopaedic Trauma"

# --- Row 3 (LeBrun et al., Hip Fracture Outcomes...) ---
$ws.Range("D3").Value = "Supplemental Digital Content is Available in the Text.`n"
$ws.Range("E3").Value = "[Drake G.%LeBrun%NULL%0, Maxwell A.%Konnaris%NULL%2, Gregory C.%Ghahramani%NULL%2, Ajay%Premkumar%NULL%2, Chris J.%DeFrancesco%NULL%2, Jordan A.%Gruskay%NULL%2, Aleksey%Dvorzhinskiy%NULL%2, Milan S.%Sandhu%NULL%2, Elan M.%Goldwyn%NULL%2, Christopher L.%Mendias%NULL%2, William M.%Ricci%NULL%2]"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "Journal of Orthopaedic Trauma"

# --- Row 4 (Muñoz Vives et al., Mortality Rates...) ---
$ws.Range("E4").Value = "[Josep Maria%Muñoz Vives%NULL%0, Montsant%Jornet-Gibert%NULL%2, Montsant%Jornet-Gibert%NULL%0, J.%Cámara-Cabrera%NULL%2, J.%Cámara-Cabrera%NULL%0, Pedro L.%Esteban%NULL%2, Pedro L.%Esteban%NULL%0, Laia%Brunet%NULL%2, Laia%Brunet%NULL%0, Luis%Delgado-Flores%NULL%2, Luis%Delgado-Flores%NULL%0, P.%Camacho-Carrasco%NULL%2, P.%Camacho-Carrasco%NULL%0, P.%Torner%NULL%2, P.%Torner%NULL%0, Francesc%Marcano-Fernández%NULL%2, Francesc%Marcano-Fernández%NULL%0, NULL%NULL%NULL%0, NULL%NULL%NULL%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Journal of Bone and Joint Surgery, Inc."

# --- Row 5 (Cheung & Forsh, Early outcomes...) ---
$ws.Range("E5").Value = "[Zoe B.%Cheung%zoe.cheung@mountsinai.org%0, David A.%Forsh%NULL%1]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Elsevier"

# --- Row 6 (IMPACT-Scot CrossRef-only record) reset to unknowns ---
$ws.Range("C6").Value = "Unknown Title"
$ws.Range("D6").Value = "Unknown Abstract"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "not found"
$ws.Range("G6").Value = "N/A"
# "1970-01-01" reads as a date literal, so Excel would otherwise convert it
# to a date serial number. Forcing the Text number format before assigning
# keeps it a literal string; ClearFormats() afterwards drops the formatting
# override again while leaving the text itself intact.
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "1970-01-01"
$ws.Range("H6").ClearFormats()
$ws.Range("J6").Value = ""

# --- Row 7 (perioperative morbidity CrossRef-only record) reset to unknowns ---
$ws.Range("C7").Value = "Unknown Title"
$ws.Range("D7").Value = "Unknown Abstract"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = "not found"
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "1970-01-01"
$ws.Range("H7").ClearFormats()
$ws.Range("J7").Value = ""

Write-Host "Applied citation graph fixes"
